$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two self-pairing rows (sending cluster == target cluster):
# old row 9 (MuSCs -> MuSCs) and old row 4 (FAPs -> MuSCs)
$ws.Rows("9").Delete()
$ws.Rows("4").Delete()

# Apply updated (re-normalised TPM) values to the remaining rows
# Row 2
$ws.Range("M2").Value = 1.991853333333333
$ws.Range("N2").Value = 5.97556
$ws.Range("O2").Value = 0.02613402671393343
$ws.Range("P2").Value = 0.02613402671393343
$ws.Range("Q2").Value = 3.548340644088888
$ws.Range("R2").Value = 31.9350657968
$ws.Range("S2").Value = 0.02571618940179495
$ws.Range("T2").Value = 0.02571618940179494
# Row 3
$ws.Range("M3").Value = 36.19366066666667
$ws.Range("N3").Value = 108.580982
$ws.Range("O3").Value = 0.4748773812350852
$ws.Range("P3").Value = 0.4748773812350851
$ws.Range("Q3").Value = 64.47635227588444
$ws.Range("R3").Value = 580.2871704829599
$ws.Range("S3").Value = 0.4672849236799377
$ws.Range("T3").Value = 0.4672849236799376
# Row 4
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.06671
$ws.Range("N4").Value = 18.20013
$ws.Range("O4").Value = 0.07959800982954926
$ws.Range("P4").Value = 0.07959800982954925
$ws.Range("Q4").Value = 10.80739897293333
$ws.Range("R4").Value = 97.26659075639999
$ws.Range("S4").Value = 0.07832537707215564
$ws.Range("T4").Value = 0.07832537707215563
# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("M5").Value = 31.96463133333333
$ws.Range("N5").Value = 95.893894
$ws.Range("O5").Value = 0.4193905822214322
$ws.Range("P5").Value = 0.4193905822214321
$ws.Range("Q5").Value = 56.94264664736888
$ws.Range("R5").Value = 512.48381982632
$ws.Range("S5").Value = 0.4126852613946891
$ws.Range("T5").Value = 0.412685261394689
# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02894466666666666
$ws.Range("H6").Value = 0.08683399999999999
$ws.Range("I6").Value = 0.01598824845142267
$ws.Range("J6").Value = 0.01598824845142267
$ws.Range("M6").Value = 1.991853333333333
$ws.Range("N6").Value = 5.97556
$ws.Range("O6").Value = 0.02613402671393343
$ws.Range("P6").Value = 0.02613402671393343
$ws.Range("Q6").Value = 0.05765353078222221
$ws.Range("R6").Value = 0.5188817770399999
$ws.Range("S6").Value = 0.000417837312138485
$ws.Range("T6").Value = 0.000417837312138485
# Row 7
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("G7").Value = 0.02894466666666666
$ws.Range("H7").Value = 0.08683399999999999
$ws.Range("I7").Value = 0.01598824845142267
$ws.Range("J7").Value = 0.01598824845142267
$ws.Range("M7").Value = 36.19366066666667
$ws.Range("N7").Value = 108.580982
$ws.Range("O7").Value = 0.4748773812350852
$ws.Range("P7").Value = 0.4748773812350851
$ws.Range("Q7").Value = 1.047613443443111
$ws.Range("R7").Value = 9.428520990988
$ws.Range("S7").Value = 0.007592457555147505
$ws.Range("T7").Value = 0.007592457555147504
# Row 8
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("G8").Value = 0.02894466666666666
$ws.Range("H8").Value = 0.08683399999999999
$ws.Range("I8").Value = 0.01598824845142267
$ws.Range("J8").Value = 0.01598824845142267
$ws.Range("M8").Value = 6.06671
$ws.Range("N8").Value = 18.20013
$ws.Range("O8").Value = 0.07959800982954926
$ws.Range("P8").Value = 0.07959800982954925
$ws.Range("Q8").Value = 0.1755988987133333
$ws.Range("R8").Value = 1.58039008842
$ws.Range("S8").Value = 0.001272632757393618
$ws.Range("T8").Value = 0.001272632757393618
# Row 9
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 0.02894466666666666
$ws.Range("H9").Value = 0.08683399999999999
$ws.Range("I9").Value = 0.01598824845142267
$ws.Range("J9").Value = 0.01598824845142267
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 31.96463133333333
$ws.Range("N9").Value = 95.893894
$ws.Range("O9").Value = 0.4193905822214322
$ws.Range("P9").Value = 0.4193905822214321
$ws.Range("Q9").Value = 0.9252055990662221
$ws.Range("R9").Value = 8.326850391596
$ws.Range("S9").Value = 0.006705320826743066
$ws.Range("T9").Value = 0.006705320826743065
